$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that follows the title.
# -----------------------------------------------------------------------
$d.Paragraphs.Item(2).Range.Delete()

# -----------------------------------------------------------------------
# 2. Split the closing "Prompt for DALLE" paragraph into two paragraphs:
#    a new bold paragraph containing the page title, followed by the
#    (still italic) paragraph whose text becomes the meta description.
# -----------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$start = $lastPara.Range.Start

$boldText = "Play Egyptian Heroes Free: Review of NetEnt's Slot Game"
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $boldText + '</w:t></w:r></w:p>'

$insertRange = $d.Range($start, $start)
$insertRange.InsertXML($xmlFrag)

$breakPos = $start + $boldText.Length
$breakRange = $d.Range($breakPos, $breakPos)
$breakRange.InsertParagraphAfter()

# -----------------------------------------------------------------------
# 3. Replace the DALLE-prompt text with the new meta-description text.
# -----------------------------------------------------------------------
$oldText = 'Prompt for DALLE: Create a feature image fitting the game "Egyptian Heroes". The image must be in a cartoon style and should feature a happy Maya warrior with glasses. The inspiration for the Maya warrior design could come from traditional Maya clothes and accessories such as the headdress and the jewelry. The background should have a modern and stylish vibe, with colors that match the superheroic connotations of the game. Some inspirational elements could be futuristic buildings, neon lights, or explosions. The text "Egyptian Heroes by NetEnt" should be included in the image. The image should convey a sense of excitement and adventure.'
$newText = "Read our review of Egyptian Heroes by NetEnt, an Egyptian themed slot game with bonus features, free spins, and an original twist. Play for free now."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
